$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cliente")

# Row 15 currently empty (style-only). Fill it like the rows above (row 14)
# by copying formatting from row 14 and then setting the new values.
$ws.Range("B14:G14").Copy()
$ws.Range("B15:G15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B15").Value = "FELIPE S."
$ws.Range("C15").Value = "d700bacde97c58b7cc8d0b53476b5697"
$ws.Range("D15").Value = Get-Date -Year 2022 -Month 10 -Day 9 -Hour 0 -Minute 0 -Second 0
$ws.Range("E15").Value = 365
$ws.Range("F15").Value = "-"
$ws.Range("G15").Value = "VENDA 11 (09/10)"
